# Update Name of Algo
# Apply updated KNN-imputed values for columns A and B on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new numeric value, taken from the target diff.
$updates = @{
    "B4"  = 7.58
    "A8"  = -21.854
    "A10" = -21.721
    "A12" = -21.53
    "B12" = 5.397
    "B15" = 5.776
    "B17" = 4.944
    "A18" = -21.885
    "B26" = 5.619
    "B27" = 5.927999999999999
    "B28" = 5.815
    "A37" = -21.032
    "B37" = 7.888000000000001
    "B47" = 5.665999999999999
    "A55" = -22.184
    "B65" = 5.912999999999999
    "A68" = -21.567
    "B73" = 6.528
    "A77" = -20.93
    "A78" = -20.582
    "A81" = -21.714
    "A82" = -21.886
    "B84" = 5.593
    "B85" = 5.964
    "B93" = 5.587000000000001
    "B95" = 5.898999999999999
    "B98" = 6.423
    "B99" = 5.415000000000001
    "B101" = 5.372
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
